# "Me agregué más tareas." — add "Matías" as Encargado for a handful of
# existing task rows (E8:E12), and turn row 14's task cell into a
# wrapped / underlined note with a matching blank "Encargado" placeholder,
# mirroring the formatting already used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Encargado" (assignee) entries for rows 8-12.
foreach ($r in 8..12) {
    $ws.Range("E$r").Value = "Matías"
}

# Row 14: task text now wraps and is underlined, so the row grows a bit
# taller, and an empty formatted "Encargado" placeholder cell appears
# next to it (matching the blank filler cells already used in other
# rows of the sheet).
$ws.Range("C14").WrapText = $true
$ws.Range("C14").Font.Underline = $true
$ws.Rows.Item(14).RowHeight = 23.95

# E14 becomes a blank "Encargado" placeholder like the H14:K14 filler
# cells already on that row — copy their format across rather than
# typing a value, so the cell stays empty but formatted.
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Restore the view to the top of the sheet with C13 selected.
$ws.Range("C13").Select() | Out-Null
